$wb = $excel.ActiveWorkbook

# --- Update selections on existing sheets (cosmetic view-state changes) ---
$wsColumnas = $wb.Worksheets.Item("columnas")
$wsColumnas.Activate()
$wsColumnas.Range("B9").Select() | Out-Null

$wsValores = $wb.Worksheets.Item("valores")
$wsValores.Activate()
$wsValores.Range("B4").Select() | Out-Null

$wsTendencia = $wb.Worksheets.Item("tendencia")
$wsTendencia.Activate()
$wsTendencia.Range("A11").Select() | Out-Null

# --- Add the new "Sheet1" worksheet after "tendencia" (becomes last/active tab) ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "Sheet1"

# Populate header row (order chosen to reproduce shared-string table ordering)
$ws.Range("A1").Value = "Elección"
$ws.Range("E1").Value = "%"
$ws.Range("F1").Value = "'-1-1"
$ws.Range("B1").Value = "Mesa"
$ws.Range("G1").Value = "pendiente1"
$ws.Range("H1").Value = "pendiente2"
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 80
$ws.Range("D2").Value = 161
$ws.Range("E2").Formula = "=C2/(C2+D2)"
$ws.Range("F2").Formula = "=D2-C2"

# Row 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 134
$ws.Range("D3").Value = 126
$ws.Range("E3").Formula = "=C3/(C3+D3)"
$ws.Range("F3").Formula = "=D3-C3"

# Row 4 (summary formulas)
$ws.Range("F4").Formula = "=F2-F3"
$ws.Range("G4").Formula = "=(D3-D2)/(C3-C2)"
$ws.Range("H4").Formula = "=(C3-C2)/(D3-D2)"

# Widen the trend-slope columns to fit their content
$ws.Range("G1:H1").ColumnWidth = 12.7109375

# Make the new sheet the active / selected tab
$ws.Activate()
$ws.Range("E17").Select() | Out-Null

Write-Output "done"
